$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 12: "Search Functionality" -> merged Search & Back test, now automated ---
$ws.Range("B12").Value = "Validate Search & Back Functionality In Search Screen"
$ws.Range("D12").Value = "Automation"
$ws.Range("E12").Value = "Automated"

# --- Row 13: "Scroll In Search Screen" -> "Scroll Up In Search Results" ---
$ws.Range("B13").Value = "Validate Scroll Up In Search Results"
$ws.Range("D13").Value = "Automation"
$ws.Range("E13").Value = "Pending"

# --- Row 14: "Back Button Functionality In Search Screen" -> "Scroll Down In Search Results" ---
$ws.Range("B14").Value = "Validate Scroll Down In Search Results"
$ws.Range("D14").Value = "Automation"
$ws.Range("E14").Value = "Pending"

# --- Row 15: View Details Of Search Record, now automated ---
$ws.Range("D15").Value = "Automation"
$ws.Range("E15").Value = "Automated"

# --- Row 16: Add To Favourites Of Search Record, now automated ---
$ws.Range("D16").Value = "Automation"
$ws.Range("E16").Value = "Automated"

# --- Row 17: Sort Search Record 'Near Me' -> 'Rating', now automated ---
$ws.Range("B17").Value = "Validate Sort Search Record Using 'Rating' Option"
$ws.Range("D17").Value = "Automation"
$ws.Range("E17").Value = "Automated"

# --- Row 18: Refresh By Dragging Down In Search Result Screen, now automated (pending) ---
$ws.Range("D18").Value = "Automation"
$ws.Range("E18").Value = "Pending"

# --- Row 19: app running in background, now automated ---
$ws.Range("D19").Value = "Automation"
$ws.Range("E19").Value = "Automated"

# --- Row 20: app in Landscape mode, now automated ---
$ws.Range("D20").Value = "Automation"
$ws.Range("E20").Value = "Automated"

# --- Row 21: brand-new test case appended at the end ---
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "Validate app when no internet connection"
$ws.Range("C21").Value = "Low"
$ws.Range("D21").Value = "Automation"
$ws.Range("E21").Value = "Automated"

# Selection moves to D5, matching the saved view state in the target file
$ws.Range("D5").Select()
